# updated steps,feature,pageobj for array,queue
#
# - Rename sheet "ArrayTryCode" -> "TryCode"
# - Make "TryCode" the active tab, with cell L20 selected
# - On "PracticeQns": scroll to row 7 (keeping the existing B8 selection),
#   no longer the active tab, and give B8 / B10 a Text number format
#   (adds a new cellXfs entry: numFmtId 49 "@")

$wb = $excel.ActiveWorkbook

$wsTry = $wb.Worksheets.Item("ArrayTryCode")
$wsTry.Name = "TryCode"

$wsQns = $wb.Worksheets.Item("PracticeQns")
$wsQns.Range("B8").NumberFormat = "@"
$wsQns.Range("B10").NumberFormat = "@"

# Scroll PracticeQns so row 7 becomes the top-left visible row, while
# keeping its real selection on B8.
$wsQns.Activate()
$excel.ActiveWindow.ScrollRow = 7
$wsQns.Range("B8").Select()

# TryCode becomes the active/selected sheet tab, with L20 selected.
$wsTry.Activate()
$wsTry.Range("L20").Select()
